$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -and $val.Contains(",")) {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $reversed = $trimmed[($trimmed.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
